$wb = $excel.ActiveWorkbook

# --- Sheet 3: rename "3) Built model with equations" -> "3)Equations" ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "3)Equations"

# --- Remove the two obsolete "GPP xxx_L rate" rows (old rows 19-20) ---
$ws3.Rows.Item(19).Delete()
$ws3.Rows.Item(19).Delete()

# --- Insert 4 blank rows for the new "NEP and Oxygen Flux" section
#     (before the old "Other" row, which is now row 23) ---
$ws3.Rows.Item(23).Resize(4).Insert()

# Fill in the new content in the same order the original author typed it in
# (this keeps the shared-string table ordering identical to the source file)
$ws3.Range("A23").Value = "NEP and Oxygen Flux"
$ws3.Range("A25").Value = "Fatm = 0.7 * (DOconc - DOsat)/Zmix"
$ws3.Range("A26").Value = "DO(t+1) = DOconc + NEP - Fatm"
$ws3.Range("A24").Value = "NEP (as O2) = (NPP - DOCrespired) * 32/12"
$ws3.Range("B24").Value = "g O2/m3"
$ws3.Range("B25").Value = "g O2/m2"
$ws3.Range("B26").Value = "g O2/m3"

# Update the NPP formulas: replace the "(1-R_auto)" factor with a flat 0.2
$ws3.Range("A19").Value = "NPP DOC_L = GPP DOC_L * 0.2 * Area / 1000"
$ws3.Range("A20").Value = "NPP POC_L = GPP POC_L * 0.2 * Area / 1000"

# Format the new section header (A23) like the other section headers
# (bold, centered, Arial 10 -- matching the "Sedimentation"/"Other" headers)
$ws3.Range("A23").Font.Bold = $true
$ws3.Range("A23").Font.Name = "Arial"
$ws3.Range("A23").Font.Size = 10
$ws3.Range("A23").HorizontalAlignment = -4108

# Give sheet3 a portrait page setup (matches the new printer/page setup entry)
$ws3.PageSetup.Orientation = 1

# --- Make "3)Equations" the active tab (previously sheet 5 "5) RMSE" was active) ---
$ws3.Activate()
$excel.ActiveWindow.Zoom = 100
$ws3.Range("A9").Select() | Out-Null
